# Normalize Portuguese accented characters to their plain-ASCII equivalents
# in the NOME_UNIDADE (column B) and NIVEL_CURSO (column C) columns.
# Example: "Graduação" -> "Graduacao"
#          "Curso Superior de Tecnologia em Agronegócio" -> "Curso Superior de Tecnologia em Agronegocio"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

function Remove-Accents($text) {
    if ($null -eq $text) { return $text }
    $result = $text
    $result = $result.Replace([char]0x00E1, 'a')  # á
    $result = $result.Replace([char]0x00E0, 'a')  # à
    $result = $result.Replace([char]0x00E2, 'a')  # â
    $result = $result.Replace([char]0x00E3, 'a')  # ã
    $result = $result.Replace([char]0x00E4, 'a')  # ä
    $result = $result.Replace([char]0x00E7, 'c')  # ç
    $result = $result.Replace([char]0x00E8, 'e')  # è
    $result = $result.Replace([char]0x00E9, 'e')  # é
    $result = $result.Replace([char]0x00EA, 'e')  # ê
    $result = $result.Replace([char]0x00EB, 'e')  # ë
    $result = $result.Replace([char]0x00EC, 'i')  # ì
    $result = $result.Replace([char]0x00ED, 'i')  # í
    $result = $result.Replace([char]0x00EE, 'i')  # î
    $result = $result.Replace([char]0x00EF, 'i')  # ï
    $result = $result.Replace([char]0x00F2, 'o')  # ò
    $result = $result.Replace([char]0x00F3, 'o')  # ó
    $result = $result.Replace([char]0x00F4, 'o')  # ô
    $result = $result.Replace([char]0x00F5, 'o')  # õ
    $result = $result.Replace([char]0x00F6, 'o')  # ö
    $result = $result.Replace([char]0x00F9, 'u')  # ù
    $result = $result.Replace([char]0x00FA, 'u')  # ú
    $result = $result.Replace([char]0x00FB, 'u')  # û
    $result = $result.Replace([char]0x00FC, 'u')  # ü
    $result = $result.Replace([char]0x00F1, 'n')  # ñ
    $result = $result.Replace([char]0x00C1, 'A')  # Á
    $result = $result.Replace([char]0x00C0, 'A')  # À
    $result = $result.Replace([char]0x00C2, 'A')  # Â
    $result = $result.Replace([char]0x00C3, 'A')  # Ã
    $result = $result.Replace([char]0x00C4, 'A')  # Ä
    $result = $result.Replace([char]0x00C7, 'C')  # Ç
    $result = $result.Replace([char]0x00C8, 'E')  # È
    $result = $result.Replace([char]0x00C9, 'E')  # É
    $result = $result.Replace([char]0x00CA, 'E')  # Ê
    $result = $result.Replace([char]0x00CB, 'E')  # Ë
    $result = $result.Replace([char]0x00CC, 'I')  # Ì
    $result = $result.Replace([char]0x00CD, 'I')  # Í
    $result = $result.Replace([char]0x00CE, 'I')  # Î
    $result = $result.Replace([char]0x00CF, 'I')  # Ï
    $result = $result.Replace([char]0x00D2, 'O')  # Ò
    $result = $result.Replace([char]0x00D3, 'O')  # Ó
    $result = $result.Replace([char]0x00D4, 'O')  # Ô
    $result = $result.Replace([char]0x00D5, 'O')  # Õ
    $result = $result.Replace([char]0x00D6, 'O')  # Ö
    $result = $result.Replace([char]0x00D9, 'U')  # Ù
    $result = $result.Replace([char]0x00DA, 'U')  # Ú
    $result = $result.Replace([char]0x00DB, 'U')  # Û
    $result = $result.Replace([char]0x00DC, 'U')  # Ü
    $result = $result.Replace([char]0x00D1, 'N')  # Ñ
    return $result
}

for ($r = 2; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    $bVal = $bCell.Value()
    $cVal = $cCell.Value()

    if ($bVal -is [string]) {
        $newB = Remove-Accents $bVal
        if ($newB -ne $bVal) {
            $bCell.Value = $newB
        }
    }

    if ($cVal -is [string]) {
        $newC = Remove-Accents $cVal
        if ($newC -ne $cVal) {
            $cCell.Value = $newC
        }
    }
}
